$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '96.971.04'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.16%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.703.69'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.29%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '237.97'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -2.71%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.93'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +1.31%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '656.44'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -1.61%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.428'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.47%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -3.56%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.06%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '3.702.32'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.39%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '44.29'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -3.48%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.47%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000296'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +10.81%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.80'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +2.99%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.395.77'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.30%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '96.709.09'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.07%  '
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.86%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.712.16'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.58%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.06'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.39%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '18.71'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +1.01%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.515'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -4.42%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '523.32'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.45%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -1.25%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000210'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.06%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.95'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.16%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '101.95'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.04%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.194'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +15.88%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '13.43'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +3.10%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '12.30'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +1.77%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.01'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.56%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.999'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.16%  '
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = 'Fetch.AI'
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.88'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +7.80%  '
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = 'Cronos'
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.189'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.31%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.02%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '32.28'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -2.29%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '649.27'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +5.88%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +1.75%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.85'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.66%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.02%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.86'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +10.75%  '
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = 'ImmutableX'
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.04'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +4.66%  '
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = 'Kaspa'
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.160'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.43%  '
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '40.38'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -5.71%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.961'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.05%  '
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0456'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.82%  '
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = 'Algorand'
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.444'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +3.69%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.88%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '23.64'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.00%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.64%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.52'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.61%  '
